$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add ages for existing rows
$ws.Range("B2").Value = 35
$ws.Range("B3").Value = 26

# Add new row with name and age
$ws.Range("A4").Value = "Amy"
$ws.Range("B4").Value = 21

# Update selection to mirror the recorded cursor position after the edit
$ws.Range("A5").Select()
